$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("E1").Value2 = "Current Difference"

# Add new header H1 ("P/L (%)"), copying the bold/border header style from G1
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H1").Value2 = "P/L (%)"

# --- Ensure columns E:H (data rows) are formatted as Text so the ---
# --- numeric-looking / percent-looking strings are stored literally ---
# --- instead of being auto-converted to numbers by the Value setter. ---
$ws.Range("E2:H26").NumberFormat = "@"

# --- Per-row data updates ---
# Row 2
$ws.Range("C2").Value2 = 182.01
$ws.Range("E2").Value2 = "61.65"
$ws.Range("F2").Value2 = "123.31"
$ws.Range("G2").Value2 = "51.23 %"
$ws.Range("H2").Value2 = "26.53 %"

# Row 3
$ws.Range("C3").Value2 = 116.14
$ws.Range("E3").Value2 = "31.54"
$ws.Range("F3").Value2 = "31.54"
$ws.Range("G3").Value2 = "37.28 %"
$ws.Range("H3").Value2 = "26.53 %"

# Row 4
$ws.Range("C4").Value2 = 168.21
$ws.Range("E4").Value2 = "88.37"
$ws.Range("F4").Value2 = "88.37"
$ws.Range("G4").Value2 = "110.68 %"
$ws.Range("H4").Value2 = "26.53 %"

# Row 5
$ws.Range("C5").Value2 = 207.87
$ws.Range("E5").Value2 = "59.92"
$ws.Range("F5").Value2 = "59.92"
$ws.Range("G5").Value2 = "40.50 %"
$ws.Range("H5").Value2 = "26.53 %"

# Row 6
$ws.Range("C6").Value2 = 46.18
$ws.Range("E6").Value2 = "17.21"
$ws.Range("F6").Value2 = "103.24"
$ws.Range("G6").Value2 = "59.39 %"
$ws.Range("H6").Value2 = "26.53 %"

# Row 7
$ws.Range("C7").Value2 = 45897.945
$ws.Range("E7").Value2 = "1064.94"
$ws.Range("F7").Value2 = "31.95"
$ws.Range("G7").Value2 = "2.38 %"
$ws.Range("H7").Value2 = "26.53 %"

# Row 8
$ws.Range("C8").Value2 = 806.5
$ws.Range("E8").Value2 = "87.50"
$ws.Range("F8").Value2 = "87.50"
$ws.Range("G8").Value2 = "12.17 %"
$ws.Range("H8").Value2 = "26.53 %"

# Row 9
$ws.Range("C9").Value2 = 63.16
$ws.Range("E9").Value2 = "22.83"
$ws.Range("F9").Value2 = "45.66"
$ws.Range("G9").Value2 = "56.61 %"
$ws.Range("H9").Value2 = "26.53 %"

# Row 10
$ws.Range("C10").Value2 = 104.15
$ws.Range("E10").Value2 = "31.67"
$ws.Range("F10").Value2 = "158.35"
$ws.Range("G10").Value2 = "43.69 %"
$ws.Range("H10").Value2 = "26.53 %"

# Row 11
$ws.Range("C11").Value2 = 156.73
$ws.Range("E11").Value2 = "55.99"
$ws.Range("F11").Value2 = "55.99"
$ws.Range("G11").Value2 = "55.58 %"
$ws.Range("H11").Value2 = "26.53 %"

# Row 12
$ws.Range("C12").Value2 = 133.85
$ws.Range("E12").Value2 = "-20.65"
$ws.Range("F12").Value2 = "-20.65"
$ws.Range("G12").Value2 = "-13.37 %"
$ws.Range("H12").Value2 = "26.53 %"

# Row 13
$ws.Range("C13").Value2 = 53.21
$ws.Range("E13").Value2 = "5.05"
$ws.Range("F13").Value2 = "20.21"
$ws.Range("G13").Value2 = "10.49 %"
$ws.Range("H13").Value2 = "26.53 %"

# Row 14
$ws.Range("C14").Value2 = 59.29499999999999
$ws.Range("E14").Value2 = "13.07"
$ws.Range("F14").Value2 = "52.27"
$ws.Range("G14").Value2 = "28.27 %"
$ws.Range("H14").Value2 = "26.53 %"

# Row 15
$ws.Range("C15").Value2 = 94.09999999999999
$ws.Range("E15").Value2 = "40.00"
$ws.Range("F15").Value2 = "80.00"
$ws.Range("G15").Value2 = "73.94 %"
$ws.Range("H15").Value2 = "26.53 %"

# Row 16
$ws.Range("C16").Value2 = 334.75
$ws.Range("E16").Value2 = "93.31"
$ws.Range("F16").Value2 = "93.31"
$ws.Range("G16").Value2 = "38.65 %"
$ws.Range("H16").Value2 = "26.53 %"

# Row 17
$ws.Range("C17").Value2 = 56.66
$ws.Range("E17").Value2 = "24.83"
$ws.Range("F17").Value2 = "24.83"
$ws.Range("G17").Value2 = "78.01 %"
$ws.Range("H17").Value2 = "26.53 %"

# Row 18
$ws.Range("C18").Value2 = 162.89
$ws.Range("E18").Value2 = "27.09"
$ws.Range("F18").Value2 = "54.18"
$ws.Range("G18").Value2 = "19.95 %"
$ws.Range("H18").Value2 = "26.53 %"

# Row 19
$ws.Range("C19").Value2 = 95.76000000000001
$ws.Range("E19").Value2 = "11.61"
$ws.Range("F19").Value2 = "23.21"
$ws.Range("G19").Value2 = "13.79 %"
$ws.Range("H19").Value2 = "26.53 %"

# Row 20
$ws.Range("C20").Value2 = 22.48
$ws.Range("E20").Value2 = "11.57"
$ws.Range("F20").Value2 = "34.71"
$ws.Range("G20").Value2 = "106.05 %"
$ws.Range("H20").Value2 = "26.53 %"

# Row 21
$ws.Range("C21").Value2 = 23.84
$ws.Range("E21").Value2 = "-7.60"
$ws.Range("F21").Value2 = "-22.79"
$ws.Range("G21").Value2 = "-24.16 %"
$ws.Range("H21").Value2 = "26.53 %"

# Row 22
$ws.Range("C22").Value2 = 126.2599
$ws.Range("E22").Value2 = "66.12"
$ws.Range("F22").Value2 = "66.12"
$ws.Range("G22").Value2 = "109.94 %"
$ws.Range("H22").Value2 = "26.53 %"

# Row 23
$ws.Range("C23").Value2 = 25.43
$ws.Range("E23").Value2 = "-4.24"
$ws.Range("F23").Value2 = "-4.24"
$ws.Range("G23").Value2 = "-14.29 %"
$ws.Range("H23").Value2 = "26.53 %"

# Row 24
$ws.Range("C24").Value2 = 178.61
$ws.Range("E24").Value2 = "10.79"
$ws.Range("F24").Value2 = "10.79"
$ws.Range("G24").Value2 = "6.43 %"
$ws.Range("H24").Value2 = "26.53 %"

# Row 25
$ws.Range("C25").Value2 = 54.06
$ws.Range("E25").Value2 = "-3.71"
$ws.Range("F25").Value2 = "-7.42"
$ws.Range("G25").Value2 = "-6.42 %"
$ws.Range("H25").Value2 = "26.53 %"

# Row 26
$ws.Range("C26").Value2 = 265.6
$ws.Range("E26").Value2 = "117.60"
$ws.Range("F26").Value2 = "117.60"
$ws.Range("G26").Value2 = "79.46 %"
$ws.Range("H26").Value2 = "26.53 %"
